# Continue working on task
# Reworks the "cancel task" flow in the eventAction sheet: cancelling a
# task now spends money (task.breakUpFee) before the cancellation actually
# goes through, and tidies up the "city has no tasks" dialog chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small in-place text tweaks (no row shifting) -------------------------

# cityHasNoTasksEvent: close the dialog back into "shop", not a dead end.
$ws.Range("D273").Value = "cityHasNoTasksDialog;shop"

# cityHasNoTasksDialog: the dialog text itself no longer also carries the
# trailing ";shop" routing (that now lives on the event above).
$ws.Range("D274").Value = "dialog_city_has_no_tasks"

# cancelTask: previously had no outgoing eventList at all; now kicks off
# the new money-caching / confirm-dialog chain.
$ws.Range("D277").Value = "closeWindow;cancelTaskCacheMoney;cancelTaskSetTempDialog;cancelTaskDialog"

# --- Insert the new rows that implement the cancel-fee flow --------------

# Two new rows land right after "cancelTask" (row 277), ahead of the
# existing "cancelTaskDialog" row (278), pushing it down to 280.
$ws.Rows("278:279").Insert()

$ws.Range("A278").Value = "cancelTaskCacheMoney"
$ws.Range("C278").Value = "setNumber"
$ws.Range("D278").Value = "money=task.breakUpFee"

$ws.Range("A279").Value = "cancelTaskSetTempDialog"
$ws.Range("C279").Value = "dialogTemp"
$ws.Range("D279").Value = "cache.money"

# "cancelTaskDialog" (now row 280) now routes through the new money-enough
# check instead of straight to "cancelTaskYes".
$ws.Range("D280").Value = "dialog_cancel_task_confirm;cancelTaskMoneyEnough;shop"

# Three more new rows land after "cancelTaskDialog" (280), ahead of the
# existing "cancelTaskYes" row (currently at 281), pushing it down to 284.
$ws.Rows("281:283").Insert()

$ws.Range("A281").Value = "cancelTaskMoneyEnough"
$ws.Range("C281").Value = "condition"
$ws.Range("D281").Value = "moneyEnough;cancelTaskProcess;moneyNotEnoughDialog"

$ws.Range("A282").Value = "cancelTaskProcess"
$ws.Range("C282").Value = "eventList"
$ws.Range("D282").Value = "cancelTaskMoneySpending;cancelTaskYes;shop"

$ws.Range("A283").Value = "cancelTaskMoneySpending"
$ws.Range("C283").Value = "dataChange"
$ws.Range("D283").Value = "money;-;cache.money"

# "cancelTaskYes" (now row 284) now hands the task off to the guild as a
# forsaken task rather than a plain cancel.
$ws.Range("D284").Value = "guild;forsakeTask"

# --- Cosmetic: keep the view roughly where the author left it ------------
$ws.Activate()
$ws.Range("D274").Select() | Out-Null
